# Move the "LOT2060 -  Tecnologia de Biopolímeros  (Requisito)" bullet line
# from the top of the "Requisitos" list to the bottom (after the
# "LOT2049 -  Genética e Biotecnologia Vegetal  (Requisito)" line).

$d = $word.ActiveDocument

# Step 1: delete the LOT2060 line (text + its line break) from the start of
# the list by searching for it and replacing the match with nothing.
$find = $d.Content.Find
$find.Execute(
    "LOT2060 -  Tecnologia de Biopolímeros  (Requisito)`v",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 2
) | Out-Null

# Step 2: re-insert the same line (text + trailing line break) at the end of
# the "Requisitos" paragraph, right after the LOT2049 line.
$reqPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$reqPara.Range.InsertAfter("LOT2060 -  Tecnologia de Biopolímeros  (Requisito)`v")
